{"js": "// Apply the \"Adjunct Reappointment\" template edits:\n//  - correct the letter date\n//  - normalize ALL-CAPS placeholder text (name, email, title, department,\n//    term length, dates) to properly-cased human text\n//  - drop the stray \"Adjunct\" that duplicated \"Adjunct Lecturer\"\n//  - fix the cc: line (name casing + \"Department Chair\")\n//  - fill in the signature \"NAME\" placeholder with the actual name\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Helper: replace the first match of `find` (exact / case-sensitive) inside\n// a given paragraph with `replacement`.\nasync function replaceInParagraph(paragraph, find, replacement) {\n  const range = paragraph.getRange();\n  const results = range.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found in paragraph: \" + find);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst paras = paragraphs.items;\n\n// 1) Letter date: 02/12/2026 -> 02/13/2026\nawait replaceInParagraph(paras[0], \"02/12/2026\", \"02/13/2026\");\n\n// 2) Addressee name block: KAZIMIR I. KARWOWSKI -> Kazimir I. Karwowski\nawait replaceInParagraph(paras[3], \"KAZIMIR I. KARWOWSKI\", \"Kazimir I. Karwowski\");\n\n// 3) Addressee email: KAZANG1234@GMAIL.COM -> kazang1234@gmail.com\nawait replaceInParagraph(paras[4], \"KAZANG1234@GMAIL.COM\", \"kazang1234@gmail.com\");\n\n// 4) Salutation: Dear Dr. KARWOWSKI, -> Dear Dr. Karwowski,\nawait replaceInParagraph(paras[7], \"KARWOWSKI\", \"Karwowski\");\n\n// 5) Body paragraph: several fixes in the same paragraph.\nconst bodyPara = paras[9];\nawait replaceInParagraph(bodyPara, \"C. FRED HIGGS III\", \"C. Fred Higgs III\");\nawait replaceInParagraph(\n  bodyPara,\n  \", I am pleased to reappoint you to the position of Adjunct \",\n  \", I am pleased to reappoint you to the position of \"\n);\nawait replaceInParagraph(bodyPara, \"ADJUNCT LECTURER\", \"Adjunct Lecturer\");\nawait replaceInParagraph(\n  bodyPara,\n  \"RICE CENTER FOR ENGINEERING LEADERSHIP\",\n  \"Rice Center for Engineering Leadership\"\n);\nawait replaceInParagraph(bodyPara, \"3 YEARS\", \"3 years\");\nawait replaceInParagraph(bodyPara, \"07/01/2024\", \"July 1, 2024\");\nawait replaceInParagraph(bodyPara, \"06/30/2027\", \"June 30, 2027\");\n\n// 6) cc: line -> correct casing + \"Department Chair\"\nawait replaceInParagraph(\n  paras[26],\n  \"C. FRED HIGGS III, Department C. Fred Higgs III\",\n  \"C. Fred Higgs III, Department Chair\"\n);\n\n// 7) Signature block placeholder \"NAME\" -> actual name (also absorbs one of\n// the tab stops that followed it).\nawait replaceInParagraph(paras[30], \"NAME\\t\", \"Kazimir I. Karwowski\");\n", "ps1": "# Apply the \"Adjunct Reappointment\" template edits:\n#  - correct the letter date\n#  - normalize ALL-CAPS placeholder text (name, email, title, department,\n#    term length, dates) to properly-cased human text\n#  - drop the stray \"Adjunct\" that duplicated \"Adjunct Lecturer\"\n#  - fix the cc: line (name casing + \"Department Chair\")\n#  - fill in the signature \"NAME\" placeholder with the actual name\n\n$d = $word.ActiveDocument\n\nfunction Replace-InParagraph($paraIndex, $findText, $replaceText) {\n    $range = $d.Paragraphs($paraIndex).Range\n    $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 1) Letter date: 02/12/2026 -> 02/13/2026\nReplace-InParagraph 1 \"02/12/2026\" \"02/13/2026\"\n\n# 2) Addressee name block: KAZIMIR I. KARWOWSKI -> Kazimir I. Karwowski\nReplace-InParagraph 4 \"KAZIMIR I. KARWOWSKI\" \"Kazimir I. Karwowski\"\n\n# 3) Addressee email: KAZANG1234@GMAIL.COM -> kazang1234@gmail.com\nReplace-InParagraph 5 \"KAZANG1234@GMAIL.COM\" \"kazang1234@gmail.com\"\n\n# 4) Salutation: Dear Dr. KARWOWSKI, -> Dear Dr. Karwowski,\nReplace-InParagraph 8 \"KARWOWSKI\" \"Karwowski\"\n\n# 5) Body paragraph: several fixes in the same paragraph.\nReplace-InParagraph 10 \"C. FRED HIGGS III\" \"C. Fred Higgs III\"\nReplace-InParagraph 10 \", I am pleased to reappoint you to the position of Adjunct \" \", I am pleased to reappoint you to the position of \"\nReplace-InParagraph 10 \"ADJUNCT LECTURER\" \"Adjunct Lecturer\"\nReplace-InParagraph 10 \"RICE CENTER FOR ENGINEERING LEADERSHIP\" \"Rice Center for Engineering Leadership\"\nReplace-InParagraph 10 \"3 YEARS\" \"3 years\"\nReplace-InParagraph 10 \"07/01/2024\" \"July 1, 2024\"\nReplace-InParagraph 10 \"06/30/2027\" \"June 30, 2027\"\n\n# 6) cc: line -> correct casing + \"Department Chair\"\nReplace-InParagraph 27 \"C. FRED HIGGS III, Department C. Fred Higgs III\" \"C. Fred Higgs III, Department Chair\"\n\n# 7) Signature block placeholder \"NAME\" -> actual name (also absorbs one of\n# the tab stops that followed it).\n$tabChar = [char]9\n$nameFind = \"NAME\" + $tabChar\nReplace-InParagraph 31 $nameFind \"Kazimir I. Karwowski\"\n"}
